$wb = $excel.ActiveWorkbook

# --- Update status text everywhere it appears ("Ready for handoff" -> "In Translation") ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Shrink the now-narrower "status" columns to match the regenerated report layout ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
